$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = 45213
$ws.Range("C14").Value = 0.79236111111111107
$ws.Range("D14").Value = 0.84166666666666667
$ws.Range("F14").Value = "Meeting minutes"

$ws.Range("B15").Value = 45216
$ws.Range("C15").Value = 0.38055555555555554
$ws.Range("D15").Value = 0.57638888888888895
$ws.Range("F15").Value = "Worked towards multiprocess implementation, researched preprocessing"

$ws.Range("B16").Value = 45216
$ws.Range("C16").Value = 0.95833333333333337
$ws.Range("D16").Value = 0.22222222222222221
$ws.Range("F16").Value = "Implemented render process"

$ws.Range("B17").Value = 45223
$ws.Range("C17").Value = 0.58333333333333337
$ws.Range("D17").Value = 0.60416666666666663
$ws.Range("F17").Value = "Weekly meeting"

$ws.Range("B18").Value = 45230
$ws.Range("C18").Value = 0.58333333333333337
$ws.Range("D18").Value = 0.60416666666666663
$ws.Range("F18").Value = "Weekly meeting"

$ws.Range("B19").Value = 45234
$ws.Range("C19").Value = 0.16666666666666666
$ws.Range("D19").Value = 0.6118055555555556
$ws.Range("F19").Value = "Fixed build system, implemented watchdog, released build to PyPi"

$ws.Range("B20").Value = 45237
$ws.Range("C20").Value = 0.083333333333333329
$ws.Range("D20").Value = 0.34166666666666662
$ws.Range("F20").Value = "Meeting minutes & Shader preprocessing work"

[void]$ws.Range("D22").Select()

Write-Host "done"
